# Sprint Review 3 - slide 3 ("Sprint 3: user story's") content placeholder.
# The backlog/user-story bullets get reshuffled: paragraph 1 (server link) and
# paragraph 7 (AR puzzle, red text + trailing line breaks) stay exactly where
# they are; paragraphs 2-6 are reordered in place.
#
# We reassign text through an intermediate placeholder first so the COM
# host's run-splitting (it keeps the longest shared prefix/suffix as a
# separate run when the new text overlaps the old one) doesn't fragment a
# paragraph's single <a:r> into several runs - each placeholder shares no
# characters with either the old or the final text, so both text writes end
# up being a clean single-run replacement while formatting (rPr/solidFill)
# is preserved.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shp = $s.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange

$placeholders = @(
    "###PLACEHOLDER_2###",
    "###PLACEHOLDER_3###",
    "###PLACEHOLDER_4###",
    "###PLACEHOLDER_5###",
    "###PLACEHOLDER_6###"
)

# New order for paragraphs 2..6 (1-based paragraph index via Paragraphs(index, 1)):
$newTexts = @(
    '(3sp) Als speler zou ik graag op de locatie "de stadsfeestzaal" een quiz spelen',
    "(5sp) Als speler wil ik mij kunnen registreren in het spel",
    "(3sp) Als speler wil ik graag de opdracht te zien krijgen als ik op de locatie ben.",
    "(2) Als speler zou ik graag een indicatie krijgen hoe ver ik van de volgende opdracht verwijderd ben.",
    "(3sp) Als speler wil ik een sessie kunnen aanmaken."
)

for ($i = 0; $i -lt $placeholders.Length; $i++) {
    $paraIndex = $i + 2
    $tr.Paragraphs($paraIndex, 1).Text = $placeholders[$i]
}

for ($i = 0; $i -lt $newTexts.Length; $i++) {
    $paraIndex = $i + 2
    $tr.Paragraphs($paraIndex, 1).Text = $newTexts[$i]
}
